# Rename three channel labels across the "联通组播" and "电信组播" sheets.
#   重庆影视     -> 重庆影视剧
#   重庆科教     -> 重庆红叶
#   重庆时尚生活 -> 重庆红岩文化

$wb = $excel.ActiveWorkbook

$wsUnicom = $wb.Worksheets.Item("联通组播")
$wsUnicom.Range("A14:A17").Value = "重庆影视剧"
$wsUnicom.Range("A18:A21").Value = "重庆红叶"
$wsUnicom.Range("A30:A33").Value = "重庆红岩文化"

$wsTelecom = $wb.Worksheets.Item("电信组播")
$wsTelecom.Range("A2:A7").Value = "重庆影视剧"
$wsTelecom.Range("A14:A16").Value = "重庆红叶"
$wsTelecom.Range("A23:A28").Value = "重庆红岩文化"
